$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DoubleDQN config row (row 4): EPISODES, LEARNING_RATE, EPSILON_END
$ws.Range("B4").Value = 7000
$ws.Range("C4").Value = 0.00001
$ws.Range("F4").Value = 0.01

# Update the active cell selection on the sheet
$ws.Range("H12").Select()
